# Restored from revision of admin on 03/31/2020 05:08:01 PM.TEST Author: admin. Type: SAVE.
# The only semantic content change in the target revision is cell C10 on the
# "Rules" worksheet: its value changes from 18 to 1 (numeric).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

$ws.Range("C10").Value = 1
